$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal text value to a cell, forcing Text number format
# so Excel does not auto-coerce numeric-looking strings (e.g. "278.28",
# "0.64%") into floating point numbers / percentage-formatted numbers.
# This mirrors how the source workbook stores these as plain text cells.
function Set-TextValue {
    param($Address, [string]$Text)
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

Set-TextValue "D2" "278.28"
Set-TextValue "E2" "0.64%"
Set-TextValue "D3" "27.24"
Set-TextValue "E3" "1.84%"
Set-TextValue "D4" "4.867"
Set-TextValue "D5" "0.06429"
Set-TextValue "E5" "1.57%"
Set-TextValue "D6" "7.018"
Set-TextValue "E6" "1.21%"
Set-TextValue "D7" "1.197"
Set-TextValue "E7" "-6.70%"
Set-TextValue "D8" "0.8867"
Set-TextValue "E8" "1.52%"
Set-TextValue "D9" "0.1542"
Set-TextValue "E9" "-0.85%"
Set-TextValue "D10" "0.05141"
Set-TextValue "E10" "1.03%"
Set-TextValue "E11" "0.18%"
Set-TextValue "D12" "0.02879"
Set-TextValue "E12" "-2.50%"
Set-TextValue "D13" "0.08972"
Set-TextValue "E13" "-0.96%"
Set-TextValue "D14" "0.001566"
Set-TextValue "E14" "-0.45%"
Set-TextValue "D15" "0.0006399"
Set-TextValue "E15" "1.13%"
Set-TextValue "D16" "0.006090"
Set-TextValue "E16" "1.10%"
Set-TextValue "D17" "3.476"
Set-TextValue "E17" "0.76%"
Set-TextValue "E18" "-0.47%"
Set-TextValue "D19" "2.239"
Set-TextValue "E19" "-1.96%"
Set-TextValue "E21" "0.51%"
Set-TextValue "D22" "3.919"
Set-TextValue "E22" "-0.38%"
Set-TextValue "E23" "10.04%"
Set-TextValue "D24" "0.04403"
Set-TextValue "E24" "0.70%"
Set-TextValue "D25" "0.001175"
Set-TextValue "E25" "0.02%"
Set-TextValue "D26" "0.003880"
Set-TextValue "E26" "-7.97%"
Set-TextValue "D28" "0.0001180"
Set-TextValue "E28" "-1.73%"
Set-TextValue "E29" "1.72%"
Set-TextValue "D40" "0.04120"
Set-TextValue "E40" "0.56%"
Set-TextValue "D41" "0.006798"
Set-TextValue "E41" "-2.81%"
Set-TextValue "E42" "0.10%"
Set-TextValue "E43" "-10.76%"
Set-TextValue "D44" "0.01167"
Set-TextValue "E44" "3.81%"
Set-TextValue "D45" "0.00005320"
Set-TextValue "E45" "0.52%"
Set-TextValue "D46" "1.559"
Set-TextValue "E46" "4.88%"

Write-Host "Applied 58 cell updates"
